# Apply schedule edits to the "Weslei" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3,4,6,7,8: column F updates ---
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "['MEC-1B-Metrologia 1', -, -, -]"
$ws.Range("F6").Value = "['MEC-1B-Metrologia 1', -, -, -]"
$ws.Range("F7").Value = "['MEC-1B-Metrologia 1', -, -, -]"
$ws.Range("F8").Value = "['MEC-1B-Metrologia 1', -, -, -]"

# --- Row 18 ---
$ws.Range("B18").Value = "[-, 'MEC-1NB-Desenho tecnico mecanico']"
$ws.Range("D18").Value = "[-, 'MEC-1NB-Metrologia 1', 'ELM-1NA-Metrologia', -]"

# --- Row 19 ---
$ws.Range("B19").Value = "[-, 'MEC-1NB-Desenho tecnico mecanico']"
$ws.Range("C19").Value = "[-, -, 'MEC-1NB-Metrologia 1', -]"
$ws.Range("D19").Value = "[-, -, 'ELM-1NA-Metrologia', 'MEC-1NB-Metrologia 1']"
$ws.Range("E19").Value = "['ELM-2NA-CAD', 'ELM-2NA-CAD']"

# --- Row 20 ---
$ws.Range("B20").Value = "[-, 'MEC-1NB-Desenho tecnico mecanico']"
$ws.Range("D20").Value = "[-, -, 'ELM-1NA-Metrologia', -]"
$ws.Range("E20").Value = "[-, 'ELM-2NA-CAD']"
$ws.Range("F20").Value = "['ELM-2NA-CAD', -]"

# --- Row 21 ---
$ws.Range("B21").Value = "-"
$ws.Range("D21").Value = "[-, -, 'ELM-1NA-Metrologia', 'MEC-1NB-Metrologia 1']"
$ws.Range("E21").Value = "-"
